$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.743.35"
$ws.Range("E2").Value = "  -2.86%  "
$ws.Range("D3").Value = "2.487.06"
$ws.Range("E3").Value = "  -5.29%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "554.39"
$ws.Range("E5").Value = "  -3.72%  "
$ws.Range("D6").Value = "147.37"
$ws.Range("E6").Value = "  -4.57%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.604"
$ws.Range("E8").Value = "  -3.22%  "
$ws.Range("D9").Value = "2.486.55"
$ws.Range("E9").Value = "  -5.18%  "
$ws.Range("E10").Value = "  -7.51%  "
$ws.Range("E11").Value = "  -5.96%  "
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("E13").Value = "  -5.31%  "
$ws.Range("D14").Value = "26.36"
$ws.Range("E14").Value = "  -6.67%  "
$ws.Range("D15").Value = "2.935.45"
$ws.Range("E15").Value = "  -5.27%  "
$ws.Range("D16").Value = "0.0000170"
$ws.Range("E16").Value = "  -7.00%  "
$ws.Range("D17").Value = "61.676.14"
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("D18").Value = "2.482.11"
$ws.Range("E18").Value = "  -5.33%  "
$ws.Range("D19").Value = "11.19"
$ws.Range("E19").Value = "  -7.39%  "
$ws.Range("D20").Value = "7.02"
$ws.Range("E20").Value = "  -7.23%  "
$ws.Range("E21").Value = "  -6.57%  "
$ws.Range("D22").Value = "323.28"
$ws.Range("E22").Value = "  -6.09%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "1.80"
$ws.Range("E24").Value = "  -4.64%  "
$ws.Range("D25").Value = "64.13"
$ws.Range("E25").Value = "  -5.46%  "
$ws.Range("D26").Value = "0.0000100"
$ws.Range("E26").Value = "  -6.54%  "
$ws.Range("D27").Value = "2.611.65"
$ws.Range("E27").Value = "  -4.91%  "
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "8.41"
$ws.Range("E30").Value = "  -8.56%  "
$ws.Range("D31").Value = "536.37"
$ws.Range("E31").Value = "  -10.75%  "
$ws.Range("D32").Value = "7.58"
$ws.Range("E32").Value = "  -4.56%  "
$ws.Range("D33").Value = "0.152"
$ws.Range("E33").Value = "  -5.28%  "
$ws.Range("E34").Value = "  -6.73%  "
$ws.Range("D35").Value = "1.60"
$ws.Range("E35").Value = "  -7.34%  "
$ws.Range("D36").Value = "5.96"
$ws.Range("E36").Value = "  -9.46%  "
$ws.Range("D37").Value = "4.93"
$ws.Range("E37").Value = "  -8.04%  "
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  -3.99%  "
$ws.Range("D40").Value = "18.58"
$ws.Range("E40").Value = "  -5.64%  "
$ws.Range("D41").Value = "148.67"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("E42").Value = "  -7.67%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "40.39"
$ws.Range("E44").Value = "  -3.00%  "
$ws.Range("D45").Value = "2.36"
$ws.Range("E45").Value = "  -5.53%  "
$ws.Range("D46").Value = "149.18"
$ws.Range("E46").Value = "  -5.82%  "
$ws.Range("D47").Value = "3.64"
$ws.Range("E47").Value = "  -6.40%  "
$ws.Range("D48").Value = "21.11"
$ws.Range("E48").Value = "  -12.49%  "
$ws.Range("E49").Value = "  -8.33%  "
$ws.Range("D50").Value = "0.599"
$ws.Range("E50").Value = "  -4.62%  "
$ws.Range("D51").Value = "0.0950"
$ws.Range("E51").Value = "  -4.77%  "
